$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "63.129.02"
$ws.Cells.Item(2, 5).Value = "  +0.43%  "
$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.601.56"
$ws.Cells.Item(3, 5).Value = "  +2.55%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "583.27"
$ws.Cells.Item(5, 5).Value = "  +2.31%  "
$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "148.51"
$ws.Cells.Item(6, 5).Value = "  +2.22%  "
$ws.Cells.Item(7, 5).Value = "  +0.03%  "
$ws.Cells.Item(9, 5).Value = "  +3.53%  "
$ws.Cells.Item(10, 5).Value = "  +3.71%  "
$ws.Cells.Item(11, 5).Value = "  +0.17%  "
$dCell = $ws.Cells.Item(12, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.354"
$ws.Cells.Item(12, 5).Value = "  +0.48%  "
$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "27.36"
$ws.Cells.Item(13, 5).Value = "  +0.84%  "
$dCell = $ws.Cells.Item(14, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.067.59"
$ws.Cells.Item(14, 5).Value = "  +2.64%  "
$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "63.042.59"
$ws.Cells.Item(15, 5).Value = "  +0.41%  "
$ws.Cells.Item(16, 5).Value = "  +4.16%  "
$dCell = $ws.Cells.Item(17, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.599.06"
$ws.Cells.Item(17, 5).Value = "  +3.13%  "
$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "11.40"
$ws.Cells.Item(18, 5).Value = "  +1.43%  "
$dCell = $ws.Cells.Item(19, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "344.85"
$ws.Cells.Item(19, 5).Value = "  +3.43%  "
$ws.Cells.Item(20, 5).Value = "  +2.88%  "
$dCell = $ws.Cells.Item(21, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.81"
$ws.Cells.Item(21, 5).Value = "  +1.29%  "
$ws.Cells.Item(22, 5).Value = "  -0.13%  "
$ws.Cells.Item(23, 5).Value = "  -0.88%  "
$ws.Cells.Item(24, 5).Value = "  +3.55%  "
$dCell = $ws.Cells.Item(25, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.723.98"
$ws.Cells.Item(25, 5).Value = "  +2.76%  "
$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.169"
$ws.Cells.Item(26, 5).Value = "  -0.09%  "
$ws.Cells.Item(27, 5).Value = "  +1.47%  "
$dCell = $ws.Cells.Item(28, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$ws.Cells.Item(28, 5).Value = "  -0.21%  "
$dCell = $ws.Cells.Item(29, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "8.41"
$ws.Cells.Item(29, 5).Value = "  +1.40%  "
$dCell = $ws.Cells.Item(30, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "7.86"
$ws.Cells.Item(30, 5).Value = "  +8.96%  "
$ws.Cells.Item(31, 5).Value = "  +0.48%  "
$ws.Cells.Item(32, 5).Value = "  +5.08%  "
$ws.Cells.Item(33, 5).Value = "  +2.91%  "
$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "469.50"
$ws.Cells.Item(34, 5).Value = "  +18.20%  "
$ws.Cells.Item(35, 5).Value = "  +6.53%  "
$dCell = $ws.Cells.Item(36, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "175.20"
$ws.Cells.Item(36, 5).Value = "  -0.98%  "
$ws.Cells.Item(37, 5).Value = "  +2.53%  "
$ws.Cells.Item(38, 5).Value = "  +0.11%  "
$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "19.22"
$ws.Cells.Item(39, 5).Value = "  +1.18%  "
$ws.Cells.Item(40, 5).Value = "  +6.66%  "
$ws.Cells.Item(42, 5).Value = "  -1.21%  "
$dCell = $ws.Cells.Item(43, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "159.35"
$ws.Cells.Item(43, 5).Value = "  +5.95%  "
$ws.Cells.Item(44, 5).Value = "  +2.64%  "
$ws.Cells.Item(45, 5).Value = "  +7.77%  "
$dCell = $ws.Cells.Item(46, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "21.11"
$ws.Cells.Item(46, 5).Value = "  +2.57%  "
$ws.Cells.Item(47, 5).Value = "  +4.11%  "
$ws.Cells.Item(48, 5).Value = "  +1.44%  "
$ws.Cells.Item(49, 5).Value = "  +0.90%  "
$dCell = $ws.Cells.Item(50, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "18.72"
$ws.Cells.Item(50, 5).Value = "  +4.07%  "
$dCell = $ws.Cells.Item(51, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.73"
$ws.Cells.Item(51, 5).Value = "  +2.94%  "
